# Thu, Apr 02, 2020  5:05:21 AM
#
# 1) The single table on the deck (slide 5) gets its table style swapped
#    from {98317C45-4AC4-4A87-BDBF-56468A2E9BC9} to
#    {B673B511-FE89-4D20-968A-D8E9D1EEC540}.
# 2) The presentation's theme ("Integral" / Red Violet) is swapped for the
#    stock "Office Theme" (Office) color palette that previously only lived
#    on the notes master's theme part. Font scheme / format scheme are
#    already identical between the two themes, so only the 12 theme colors
#    actually need to change.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{B673B511-FE89-4D20-968A-D8E9D1EEC540}")

# --- 2. Theme colours -------------------------------------------------------
# Master.ColorScheme indexes map 1:1 onto the DrawingML <a:clrScheme> slots:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1 .. 10 accent6, 11 hlink, 12 folHlink
# (COM RGB() helper isn't available in this host, so the 0x00BBGGRR values
# are precomputed from the target hex colours below.)
$master = $p.SlideMaster
$scheme = $master.ColorScheme

$scheme.Colors(1).RGB  = 0          # dk1      000000
$scheme.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$scheme.Colors(3).RGB  = 6968388    # dk2      44546A
$scheme.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$scheme.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$scheme.Colors(6).RGB  = 3243501    # accent2  ED7D31
$scheme.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$scheme.Colors(8).RGB  = 49407      # accent4  FFC000
$scheme.Colors(9).RGB  = 12874308   # accent5  4472C4
$scheme.Colors(10).RGB = 4697456    # accent6  70AD47
$scheme.Colors(11).RGB = 12673797   # hlink    0563C1
$scheme.Colors(12).RGB = 7491477    # folHlink 954F72
